# Insert a new "min sample length" column just before "Sample ID" (old column
# O), shifting the existing O:BB content one column to the right (to P:BC).
#
# A real EntireColumn.Insert() also re-indexes the sheet's column-width bands
# (<cols>), which the target workbook does NOT show moving - so instead we
# use a Range.Copy into the shifted destination (this carries both values and
# styles/number formats without touching column width metadata), and then
# populate the freed-up column O with the new field.
#
# Also rename "bucket" -> "bucket_sediment" for the first two data rows'
# measurement_device (column J), per the commit message (bucket_sediment /
# bucket_chamber are now the valid options).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy cell-by-cell, right-to-left, so each cell keeps its own style/number
# format intact (a single bulk Range.Copy over the whole block collapses the
# few cells that use the alternate bordered header style, s="2", onto the
# more common s="1" - copying one cell at a time avoids that).
for ($r = 1; $r -le 7; $r++) {
    for ($c = 54; $c -ge 15; $c--) {
        $ws.Cells.Item($r, $c).Copy($ws.Cells.Item($r, $c + 1))
    }
}

# New column header: borrow the neighbouring header's style (bold/wrap, no
# border) via Copy, then overwrite its text.
$ws.Range("N1").Copy($ws.Range("O1"))
$ws.Range("O1").Value = "min sample length"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 15).Value = 45
}

# bucket -> bucket_sediment for rows 2 and 3 (measurement_device, column J)
$ws.Range("J2").Value = "bucket_sediment"
$ws.Range("J3").Value = "bucket_sediment"

# Update selection to match the recorded user action.
$ws.Range("J4").Select()
